$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 614.85
$ws.Range("F31").Value = 14
$ws.Range("G31").Value = 430.36
$ws.Range("F32").Value = 24
$ws.Range("G32").Value = 614.64
$ws.Range("F36").Value = 9
$ws.Range("G36").Value = 276.66
$ws.Range("F37").Value = 20
$ws.Range("G37").Value = 922
$ws.Range("B41").Value = 22171.16
$ws.Range("F47").Value = 27
$ws.Range("G47").Value = 983.0700000000001
$ws.Range("F50").Value = 170
$ws.Range("G50").Value = 32791.3
$ws.Range("F54").Value = 81
$ws.Range("G54").Value = 1841.94
$ws.Range("F63").Value = 258
$ws.Range("G63").Value = 24133.32
$ws.Range("F76").Value = 62
$ws.Range("G76").Value = 2074.52
$ws.Range("B77").Value = 129458.74
$ws.Range("F101").Value = 6
$ws.Range("G101").Value = 2452.02
$ws.Range("B105").Value = 8472.98
$ws.Range("F113").Value = 171
$ws.Range("G113").Value = 7677.9
$ws.Range("B116").Value = 77012.83
$ws.Range("F157").Value = 59
$ws.Range("G157").Value = 2919.32
$ws.Range("F158").Value = 65
$ws.Range("G158").Value = 3216.2
$ws.Range("B168").Value = 29485.43
$ws.Range("F171").Value = 10
$ws.Range("G171").Value = 653
$ws.Range("F175").Value = 6
$ws.Range("G175").Value = 583.2
$ws.Range("F176").Value = 29
$ws.Range("G176").Value = 3669.37
$ws.Range("B183").Value = 30661.57
$ws.Range("F187").Value = 2
$ws.Range("G187").Value = 1442.52
$ws.Range("B188").Value = 19659.66
$ws.Range("F205").Value = 43
$ws.Range("G205").Value = 3365.18
$ws.Range("B210").Value = 30390.29
$ws.Range("F214").Value = 37
$ws.Range("G214").Value = 1577.68
$ws.Range("B218").Value = 8408.959999999999
$ws.Range("F227").Value = 40
$ws.Range("G227").Value = 1859.6
$ws.Range("F231").Value = 22
$ws.Range("G231").Value = 154.88
$ws.Range("F236").Value = 27
$ws.Range("G236").Value = 1131.57
$ws.Range("B237").Value = 14345.1
$ws.Range("F247").Value = 44
$ws.Range("G247").Value = 935
$ws.Range("B248").Value = 1046
$ws.Range("F251").Value = 12
$ws.Range("G251").Value = 3191.52
$ws.Range("B264").Value = 91819.41
$ws.Range("F287").Value = 47
$ws.Range("G287").Value = 6370.38
$ws.Range("B336").Value = 254361.7
$ws.Range("F340").Value = 15
$ws.Range("G340").Value = 1105.65
$ws.Range("F342").Value = 11
$ws.Range("G342").Value = 1450.13
$ws.Range("F351").Value = 79
$ws.Range("G351").Value = 1826.48
$ws.Range("F355").Value = 14
$ws.Range("G355").Value = 1776.04
$ws.Range("F361").Value = 173
$ws.Range("G361").Value = 7124.14
$ws.Range("F363").Value = 317
$ws.Range("G363").Value = 12600.75
$ws.Range("B364").Value = 63523.91
$ws.Range("F366").Value = 8
$ws.Range("G366").Value = 1227.28
$ws.Range("B370").Value = 5035.91
$ws.Range("F392").Value = 240
$ws.Range("G392").Value = 4893.6
$ws.Range("B397").Value = 27134.36
$ws.Range("F428").Value = 115
$ws.Range("G428").Value = 3141.8
$ws.Range("F433").Value = 627
$ws.Range("G433").Value = 60568.2
$ws.Range("F436").Value = 117
$ws.Range("G436").Value = 3146.13
$ws.Range("B439").Value = 132061.03
$ws.Range("F454").Value = 27
$ws.Range("G454").Value = 336.69
$ws.Range("B459").Value = 12013.75
$ws.Range("F461").Value = 462
$ws.Range("G461").Value = 5959.8
$ws.Range("F464").Value = 656
$ws.Range("G464").Value = 8462.4
$ws.Range("F465").Value = 174
$ws.Range("G465").Value = 6088.26
$ws.Range("F469").Value = 111
$ws.Range("G469").Value = 2132.31
$ws.Range("B470").Value = 40268.17
$ws.Range("F501").Value = 24
$ws.Range("G501").Value = 1201.92
$ws.Range("F503").Value = 69
$ws.Range("G503").Value = 4253.85
$ws.Range("F505").Value = 22
$ws.Range("G505").Value = 709.9400000000001
$ws.Range("F506").Value = 9
$ws.Range("G506").Value = 242.01
$ws.Range("F513").Value = 19
$ws.Range("G513").Value = 918.84
$ws.Range("B515").Value = 39508.64
$ws.Range("F524").Value = 121
$ws.Range("G524").Value = 7347.12
$ws.Range("F529").Value = 152
$ws.Range("G529").Value = 3724
$ws.Range("F530").Value = 51
$ws.Range("G530").Value = 1249.5
$ws.Range("B539").Value = 120503.93
$ws.Range("F580").Value = 25
$ws.Range("G580").Value = 1731.75
$ws.Range("B581").Value = 18405.41
$ws.Range("F583").Value = 22
$ws.Range("G583").Value = 1743.06
$ws.Range("F586").Value = 52
$ws.Range("G586").Value = 13836.68
$ws.Range("F587").Value = 105
$ws.Range("G587").Value = 2738.4
$ws.Range("F598").Value = 11
$ws.Range("G598").Value = 653.62
$ws.Range("B600").Value = 71510.89999999999
$ws.Range("F656").Value = 2
$ws.Range("G656").Value = 54.4
$ws.Range("F657").Value = 28
$ws.Range("G657").Value = 761.6
$ws.Range("F658").Value = 16
$ws.Range("G658").Value = 435.2
$ws.Range("B659").Value = 10377.72
$ws.Range("F670").Value = 0
$ws.Range("G670").Value = 0
$ws.Range("B673").Value = 1502.6
$ws.Range("F718").Value = 85
$ws.Range("G718").Value = 7454.5
$ws.Range("B722").Value = 28471.76
$ws.Range("F751").Value = 49
$ws.Range("G751").Value = 10578.61
$ws.Range("F758").Value = 97
$ws.Range("G758").Value = 5211.81
$ws.Range("F759").Value = 184
$ws.Range("G759").Value = 27807.92
$ws.Range("B767").Value = 570473.47
$ws.Range("F788").Value = 334
$ws.Range("G788").Value = 34371.94
$ws.Range("B796").Value = 168229.48
$ws.Range("B855").Value = 5444114.09
$ws.Range("B856").Value = 5444114.09
